# Add 10 new medication names to the dictionary, then re-sort the list
# (A2:A159) alphabetically, matching the workbook's existing layout, and
# extend the formatted-but-empty row range below the data down to row 199.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Append the new entries right after the current last data row (149).
#    Writing these creates the 10 new shared-string entries in this exact
#    order.
$newWords = @(
    "투제오",
    "세비카",
    "베타미가",
    "토비애즈서방정",
    "인벨라",
    "인사돌",
    "펠루비",
    "아모잘틴",
    "크레스토",
    "칸대암로정"
)

$startRow = 150
for ($i = 0; $i -lt $newWords.Count; $i++) {
    $row = $startRow + $i
    $ws.Range("A$row").Value2 = $newWords[$i]
}

# 2) Re-sort the whole list (A2:A159) ascending, same as the sheet's
#    existing autofilter sort condition on column A.
$sortRange = $ws.Range("A2:A159")
$sortRange.Sort($ws.Range("A1"))

# 3) The sheet keeps a block of blank, row-height-formatted rows below the
#    data; extend that formatting down through row 199.
for ($r = 160; $r -le 199; $r++) {
    $ws.Rows.Item($r).RowHeight = 18
}

Write-Host "done"
